$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '62.211.96'
$ws.Range("E2").Value = '  -2.12%  '

$ws.Range("D3").Value = '2.429.49'
$ws.Range("E3").Value = '  -1.80%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.28%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '561.97'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -2.60%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '142.77'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -3.87%  '

$ws.Range("E7").Value = '  +0.15%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.527'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -2.66%  '

$ws.Range("D9").Value = '2.425.20'
$ws.Range("E9").Value = '  -2.31%  '

$ws.Range("E11").Value = '  +0.92%  '

$ws.Range("B12").Value = 'Cardano'
$ws.Range("C12").Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.351'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -3.46%  '

$ws.Range("B13").Value = 'Toncoin'
$ws.Range("C13").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.16'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -3.45%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '26.37'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -3.51%  '

$ws.Range("E15").Value = '  -6.79%  '

$ws.Range("D16").Value = '2.865.36'
$ws.Range("E16").Value = '  -2.84%  '

$ws.Range("D17").Value = '62.255.81'
$ws.Range("E17").Value = '  -1.77%  '

$ws.Range("D18").Value = '2.428.47'
$ws.Range("E18").Value = '  -2.13%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '11.02'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -4.85%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.10'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -2.81%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '323.30'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -1.82%  '

$ws.Range("E22").Value = '  -3.26%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '1.99'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +3.47%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '64.84'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -3.94%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '615.00'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -2.80%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '8.94'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +2.02%  '

$ws.Range("D28").Value = '2.550.35'
$ws.Range("E28").Value = '  -1.75%  '

$ws.Range("E29").Value = '  +0.48%  '

$ws.Range("D30").Value = '0.0₃0949'
$ws.Range("E30").Value = '  -10.45%  '

$ws.Range("E31").Value = '  -6.38%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '8.00'
$ws.Range("D32").Style = "Normal"

$ws.Range("E33").Value = '  -3.81%  '

$ws.Range("E34").Value = '  -7.47%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '4.99'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -4.17%  '

$ws.Range("E36").Value = '  +0.20%  '

$ws.Range("E37").Value = '  -7.81%  '

$ws.Range("E38").Value = '  -3.41%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '18.62'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -1.97%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '146.64'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -0.37%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '5.19'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -6.47%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.73'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -7.17%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '42.58'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +1.77%  '

$ws.Range("E44").Value = '  +0.02%  '

$ws.Range("E45").Value = '  -9.39%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '145.10'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -3.87%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.67'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -2.82%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0520'
$ws.Range("D48").Style = "Normal"

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '19.97'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -6.18%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.591'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -3.21%  '

$ws.Range("E51").Value = '  -5.60%  '
